# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets to the
# freshly scraped counts (output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 3..43) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 5214
    4  = 11
    5  = 7519
    7  = 73
    9  = 604
    11 = 33
    12 = 4349
    14 = 110
    16 = 2938
    19 = 214
    20 = 518
    21 = 448
    23 = 318
    25 = 1702
    26 = 1202
    28 = 1394
    34 = 3
    35 = 63
    37 = 70
    38 = 2956
    40 = 29
    41 = 98
    43 = 55
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" (rows 3..45) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 5214
    4  = 11
    5  = 7519
    7  = 73
    9  = 604
    11 = 33
    12 = 4349
    14 = 110
    16 = 2938
    19 = 214
    20 = 518
    21 = 448
    24 = 318
    26 = 1702
    27 = 1202
    29 = 1394
    35 = 3
    36 = 63
    38 = 70
    39 = 2956
    42 = 29
    43 = 98
    45 = 55
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
